# Updates cryptos list price/volume data per the Thu Jul 18 03:31:55 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking numeric strings as TEXT (the site
# uses "." as a thousands separator, e.g. "64.775.03", which is not a valid
# number). Force a text number format first so Excel does not auto-convert
# the few cells whose new value happens to parse as a genuine number.
$textCells = @("D5", "D6", "D7", "D8", "D10", "D12", "D16", "D19", "D20", "D21", "D22", "D25", "D28", "D29", "D31", "D33", "D34", "D35", "D36", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Coin / Link / Price / Volume(1h) values.
$ws.Range('D2').Value = '64.775.03'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.426.14'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '573.98'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = '159.10'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').Value = '  +3.49%  '
$ws.Range('D9').Value = '3.426.55'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').Value = '7.17'
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').Value = '0.440'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').Value = '4.016.71'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('D16').Value = '27.72'
$ws.Range('E16').Value = '  -3.39%  '
$ws.Range('D17').Value = '64.780.64'
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '3.424.98'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '6.35'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('D20').Value = '13.91'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').Value = '380.47'
$ws.Range('E21').Value = '  -2.46%  '
$ws.Range('D22').Value = '8.02'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '72.33'
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('E27').Value = '  +6.01%  '
$ws.Range('D28').Value = '0.178'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +3.97%  '
$ws.Range('D31').Value = '6.21'
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').Value = '23.27'
$ws.Range('E33').Value = '  -1.72%  '
$ws.Range('D34').Value = '7.09'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  +4.70%  '
$ws.Range('D36').Value = '160.15'
$ws.Range('E36').Value = '  -1.46%  '
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('D38').Value = '0.0757'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').Value = '2.907.76'
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('D40').Value = '6.72'
$ws.Range('E40').Value = '  +4.34%  '
$ws.Range('D41').Value = '26.48'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('D42').Value = '4.59'
$ws.Range('E42').Value = '  +1.92%  '
$ws.Range('D43').Value = '43.01'
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('D45').Value = '0.771'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('D46').Value = '25.88'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.27'
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '317.80'
$ws.Range('E48').Value = '  +2.46%  '
$ws.Range('E49').Value = '  -3.42%  '
$ws.Range('D50').Value = '0.107'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('E51').Value = '  -1.84%  '
